$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Cell B11 currently holds the text "R40". It needs to become the text "1"
# (stored as a shared string, not a number). Simply assigning a numeric-looking
# string via .Value causes Excel to auto-detect it as a number, so instead we
# put a text-producing formula in the cell and then convert it to a static
# value via Copy / PasteSpecial(Values). This keeps the cell's existing style
# and yields a genuine text cell, matching a normal Excel edit.
$target = $ws.Range("B11")
$target.Formula = "=""1"""
$target.Copy()
$target.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false
